$p = $ppt.ActivePresentation
Write-Output "Designs.Count: $($p.Designs.Count)"

try {
    $d2 = $p.Designs.Add("Office Theme")
    Write-Output "Added, count=$($p.Designs.Count) name=$($d2.Name)"
} catch {
    Write-Output "ERR1: $_"
}
